$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the last used row in column A (component list starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Rename "Buzzer" -> "Buzzer 5v 12mm" (quantity stays the same)
$buzzerCell = $ws.Range("A2:A" + $lastRow).Find("Buzzer")
if ($buzzerCell -ne $null) {
    $buzzerCell.Value = "Buzzer 5v 12mm"
}

# Rename "Led" -> "Led Amarelo" (quantity stays the same)
$ledCell = $ws.Range("A2:A" + $lastRow).Find("Led")
if ($ledCell -ne $null) {
    $ledCell.Value = "Led Amarelo"
}

# Add a new component row: "Led Vermelho" with quantity 4.
# Seed it from the last existing data row so it inherits the same style (s="1")
# before filling in its own values.
$newRow = $lastRow + 1
$ws.Range("A" + $lastRow + ":B" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":B" + $newRow).PasteSpecial(-4122)
$ws.Range("A" + $newRow).Value = "Led Vermelho"
$ws.Range("B" + $newRow).Value = 4

# Re-sort the component table (A2:B..) alphabetically by component name,
# same as the existing sorted list.
$dataRange = $ws.Range("A2:B" + $newRow)
$keyRange = $ws.Range("A2:A" + $newRow)
$dataRange.Sort($keyRange)
